# Update cryptos list with latest prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.446.39"
$ws.Range("E2").Value = "  +12.81%  "
$ws.Range("D3").Value = "1.825.30"
$ws.Range("E3").Value = "  +9.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.77%  "
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.67"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.73%  "
$ws.Range("E10").Value = "  +6.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0672"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0931"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("D13").Value = "2.086.82"
$ws.Range("D14").Value = "1.824.22"
$ws.Range("E14").Value = "  +9.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.650"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.84%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.350.27"
$ws.Range("E16").Value = "  +12.36%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "258.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.62%  "
$ws.Range("D21").Value = "0.0₃0754"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.65%  "
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("E28").Value = "  +6.77%  "
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0522"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.01%  "
$ws.Range("E34").Value = "  +8.94%  "
$ws.Range("D35").Value = "1.553.60"
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.08%  "
$ws.Range("B38").Value = "MinaProtocolToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +216.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.637"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.52%  "
$ws.Range("E40").Value = "  +6.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "84.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.919"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.54%  "
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +9.97%  "
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("D48").Value = "1.986.78"
$ws.Range("E48").Value = "  +9.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +28.00%  "
$ws.Range("E50").Value = "  +5.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.01%  "
